$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to Text so the literal "DD-MMM-YY" strings are
# stored as-is (shared strings) instead of being auto-converted to date
# serial numbers when assigned below.
$ws.Range("A2:A5").NumberFormat = "@"

# Update date column (A) shared-string values
$ws.Range("A2").Value = "13-FEB-26"
$ws.Range("A3").Value = "20-FEB-26"
$ws.Range("A4").Value = "13-MAR-26"
$ws.Range("A5").Value = "27-MAR-26"

# Update fare threat / our fare / fare diff figures
$ws.Range("D2").Value = 687
$ws.Range("E2").Value = 895
$ws.Range("F2").Value = -208

$ws.Range("D3").Value = 513
$ws.Range("E3").Value = 883
$ws.Range("F3").Value = -370

$ws.Range("D4").Value = 1237
$ws.Range("E4").Value = 1501
$ws.Range("F4").Value = -264

$ws.Range("D5").Value = 513
$ws.Range("E5").Value = 786
$ws.Range("F5").Value = -273
